$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The Sonar mux-select labels (previously crammed into row 7 alongside the
# I2C pins) move down to row 16, next to pin "K". The stale duplicate label
# "GPIO (Sonar analog/power mux enable)" in D7 is dropped in the process.
$ws.Range("D7").ClearContents()
$ws.Range("H7:J7").ClearContents()

$ws.Range("C16").Value = "GPIO (Sonar mux select 0)"
$ws.Range("D16").Value = "GPIO (Sonar mux select 1)"
$ws.Range("E16").Value = "GPIO (Sonar mux select 2)"

# Column widths re-flow now that the long labels moved out of columns C-H.
$ws.Columns.Item(3).ColumnWidth = 30.08
$ws.Columns.Item(4).ColumnWidth = 34.25
$ws.Columns.Item(5).ColumnWidth = 29.21
$ws.Columns.Item(6).ColumnWidth = 16.52
$ws.Columns.Item(7).ColumnWidth = 31.99
$ws.Columns.Item(8).ColumnWidth = 34.6

# Restore the view to the top-left corner and match the reviewer's final
# selection / zoom level.
$ws.Range("H7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75
